$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in header: "Hyperparamter Configuration" -> "Hyperparameter Configuration"
$ws.Range("B1").Value = "Hyperparameter Configuration"

# Select the full table range (matches the saved selection state)
$ws.Range("A1:C17").Select() | Out-Null

# Adjust row heights: header row taller, data rows slightly taller
$ws.Rows("1").RowHeight = 34
$ws.Rows("2:17").RowHeight = 16

# Widen column B to fit the (now longer) header text
$ws.Columns("B").ColumnWidth = 26
